$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Cxcl1->Xcr1 TPM results: MuSCs is now also included as a
# target cluster (previously only Inflammatory-Mac and Resolving-Mac were
# targets), and all the expression/specificity metrics were recomputed with
# the new TPM values. Columns: A Sending cluster, B Ligand symbol,
# C Receptor symbol, D Target cluster, E..T numeric metrics.
$data = @"
"ECs","Cxcl1","Xcr1","Inflammatory-Mac",3,1,33.64214033333334,100.926421,0.106995191696894,0.106995191696894,3,1,1.429146333333333,4.287439,0.82551891942994,0.8255189194299399,48.07954150286879,432.715873525819,0.08832655503381923,0.0883265550338192
"ECs","Cxcl1","Xcr1","MuSCs",3,1,33.64214033333334,100.926421,0.106995191696894,0.106995191696894,1,0.3333333333333333,0.03271033333333333,0.098131,0.01889449554444493,0.01889449554444493,1.100445624350111,9.904010619151,0.002021620172793995,0.002021620172793994
"ECs","Cxcl1","Xcr1","Resolving-Mac",3,1,33.64214033333334,100.926421,0.106995191696894,0.106995191696894,2,0.6666666666666666,0.269353,0.8080590000000001,0.155586585025615,0.155586585025615,9.061611425204335,81.55450282683901,0.01664701649028078,0.01664701649028077
"FAPs","Cxcl1","Xcr1","Inflammatory-Mac",3,1,242.358393,727.0751789999999,0.7707946777896593,0.7707946777896592,3,1,1.429146333333333,4.287439,0.82551891942994,0.8255189194299399,346.365608708509,3117.290478376581,0.6363055895112684,0.6363055895112681
"FAPs","Cxcl1","Xcr1","MuSCs",3,1,242.358393,727.0751789999999,0.7707946777896593,0.7707946777896592,1,0.3333333333333333,0.03271033333333333,0.098131,0.01889449554444493,0.01889449554444493,7.927623821160999,71.34861439044899,0.01456377660517859,0.01456377660517858
"FAPs","Cxcl1","Xcr1","Resolving-Mac",3,1,242.358393,727.0751789999999,0.7707946777896593,0.7707946777896592,2,0.6666666666666666,0.269353,0.8080590000000001,0.155586585025615,0.155586585025615,65.279960229729,587.519642067561,0.1199253116732124,0.1199253116732124
"Inflammatory-Mac","Cxcl1","Xcr1","Inflammatory-Mac",3,1,9.788214000000002,29.364642,0.0311303568359039,0.03113035683590389,3,1,1.429146333333333,4.287439,0.82551891942994,0.8255189194299399,13.988790147982,125.899111331838,0.02569869853664384,0.02569869853664383
"Inflammatory-Mac","Cxcl1","Xcr1","MuSCs",3,1,9.788214000000002,29.364642,0.0311303568359039,0.03113035683590389,1,0.3333333333333333,0.03271033333333333,0.098131,0.01889449554444493,0.01889449554444493,0.3201757426780001,2.881581684102,0.0005881923885329672,0.0005881923885329669
"Inflammatory-Mac","Cxcl1","Xcr1","Resolving-Mac",3,1,9.788214000000002,29.364642,0.0311303568359039,0.03113035683590389,2,0.6666666666666666,0.269353,0.8080590000000001,0.155586585025615,0.155586585025615,2.636484805542,23.728363249878,0.004843465910727099,0.004843465910727097
"MuSCs","Cxcl1","Xcr1","Inflammatory-Mac",3,1,12.12016933333333,36.360508,0.03854688876420623,0.03854688876420622,3,1,1.429146333333333,4.287439,0.82551891942994,0.8255189194299399,17.32149556211245,155.893460059012,0.03182118596001363,0.03182118596001361
"MuSCs","Cxcl1","Xcr1","MuSCs",3,1,12.12016933333333,36.360508,0.03854688876420623,0.03854688876420622,1,0.3333333333333333,0.03271033333333333,0.098131,0.01889449554444493,0.01889449554444493,0.3964547789497778,3.568093010548,0.0007283240180075092,0.0007283240180075088
"MuSCs","Cxcl1","Xcr1","Resolving-Mac",3,1,12.12016933333333,36.360508,0.03854688876420623,0.03854688876420622,2,0.6666666666666666,0.269353,0.8080590000000001,0.155586585025615,0.155586585025615,3.264603970441334,29.381435733972,0.005997378786185098,0.005997378786185096
"Resolving-Mac","Cxcl1","Xcr1","Inflammatory-Mac",3,1,16.51773933333333,49.553218,0.0525328849133368,0.05253288491333678,3,1,1.429146333333333,4.287439,0.82551891942994,0.8255189194299399,23.60626660318912,212.456399428702,0.04336689038819519,0.04336689038819518
"Resolving-Mac","Cxcl1","Xcr1","MuSCs",3,1,16.51773933333333,49.553218,0.0525328849133368,0.05253288491333678,1,0.3333333333333333,0.03271033333333333,0.098131,0.01889449554444493,0.01889449554444493,0.5403007595064445,4.862706835558,0.0009925823599318806,0.0009925823599318802
"Resolving-Mac","Cxcl1","Xcr1","Resolving-Mac",3,1,16.51773933333333,49.553218,0.0525328849133368,0.05253288491333678,2,0.6666666666666666,0.269353,0.8080590000000001,0.155586585025615,0.155586585025615,4.449102642651334,40.041923783862,0.008173412165209726,0.008173412165209722
"@

$lines = $data -split "`n"
$row = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $fields = $line -split ","
    for ($col = 1; $col -le $fields.Count; $col++) {
        $field = $fields[$col - 1].Trim()
        if ($field.StartsWith('"') -and $field.EndsWith('"')) {
            $value = $field.Trim('"')
        } else {
            $value = [double]$field
        }
        $ws.Cells.Item($row, $col).Value = $value
    }
    $row++
}
